# Update the date in the title, and the 26 multiplication problems/answers
# in the table to their new values, per the commit "Update master to
# output generated at 503736d".
#
# Note: pair 16 (42x37=1554 -> 33x87=2871) must be applied before pair 24
# (55x47=2585 -> 42x37=1554) so that the freshly-written "42x37=1554" text
# from pair 24 is not re-matched by pair 16's search term. Processing the
# pairs in document order (top-to-bottom, as below) guarantees this.

$d = $word.ActiveDocument

$replacements = @(
    @("2025-03-31 Monday", "2025-04-01 Tuesday"),
    @("46×49=2254", "97×35=3395"),
    @("44×18=792", "20×33=660"),
    @("57×57=3249", "25×14=350"),
    @("20×18=360", "17×58=986"),
    @("22×24=528", "68×33=2244"),
    @("96×26=2496", "72×80=5760"),
    @("46×37=1702", "25×49=1225"),
    @("28×95=2660", "31×68=2108"),
    @("29×72=2088", "70×47=3290"),
    @("61×54=3294", "85×96=8160"),
    @("89×69=6141", "94×62=5828"),
    @("48×48=2304", "50×44=2200"),
    @("56×31=1736", "18×94=1692"),
    @("19×80=1520", "86×74=6364"),
    @("43×58=2494", "66×34=2244"),
    @("42×37=1554", "33×87=2871"),
    @("85×21=1785", "87×11=957"),
    @("30×42=1260", "22×36=792"),
    @("37×77=2849", "35×39=1365"),
    @("18×93=1674", "40×73=2920"),
    @("82×61=5002", "79×36=2844"),
    @("34×55=1870", "68×75=5100"),
    @("82×60=4920", "98×43=4214"),
    @("55×47=2585", "42×37=1554"),
    @("73×61=4453", "69×60=4140")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $found = $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, `
                                      $true, 1, $false, $new, 2)
    if (-not $found) {
        Write-Host "WARNING: could not find text to replace: $old"
    }
}

Write-Host "Replacements complete."
